# Fixed the name of algorithms
# Slide 7 ("Pattern Matching"), content placeholder shape (Shapes.Item(2)).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Paragraph 1: "Matching Algorithms" -> "Detecting Novelty " + "Algorithms" ---
$para1 = $tr.Paragraphs(1)
$full1 = $para1.Characters(1, $para1.Length)
$full1.Text = "Detecting Novelty Algorithms"
# Split into two runs at the word boundary so the trailing run is just "Algorithms".
$para1.Characters(1, 18).Text = "Detecting Novelty "

# --- Paragraph 2: "Statistical(Geometric) Measurement" -> "Statistical(Geometric) " + "Algorithms" ---
$para2 = $tr.Paragraphs(2)
$para2.Text = "Statistical(Geometric) Algorithms"

# --- Paragraph 3: merge the 4 runs into a single run ---
$para3 = $tr.Paragraphs(3)
$full3 = $para3.Characters(1, $para3.Length)
$full3.Text = "Nearest Neighbor using ball tree algorithm"

# --- Paragraph 4: "...Distance from the nearest neighbor" -> split into 3 runs ---
$para4 = $tr.Paragraphs(4)
# The first two runs (" " and "     : ") are untouched; only the trailing run changes.
$lastRun4 = $para4.Characters(9, 34)
$lastRun4.Text = "Distance from a point to the nearest neighbor"
$para4.Characters(9, 14).Text = "Distance from "
$para4.Characters(23, 11).Text = "a point to "
